$wb = $excel.ActiveWorkbook

# --- Metadata sheet updates ---
$meta = $wb.Worksheets.Item(1)

# Version: 0.1.8 -> 0.1.6
$meta.Cells.Item(3, 2).Value = "0.1.6"

# Status: draft -> active
$meta.Cells.Item(6, 2).Value = "active"

# Date
$meta.Cells.Item(8, 2).Value = "2023-05-05T10:50:04-05:00"

# Contact rows (both collapse to the same new text)
$meta.Cells.Item(10, 2).Value = "No display for ContactDetail"
$meta.Cells.Item(11, 2).Value = "No display for ContactDetail"

# Remove the "Jurisdiction" row entirely (old row 12), shifting later rows up
$meta.Rows.Item(12).Delete()

# --- Rename the "Include ValueSet #N" sheets to "Include ValueSets[ N]" ---
$wb.Worksheets.Item(2).Name = "Include ValueSets"
$wb.Worksheets.Item(3).Name = "Include ValueSets 2"
$wb.Worksheets.Item(4).Name = "Include ValueSets 3"
$wb.Worksheets.Item(5).Name = "Include ValueSets 4"
$wb.Worksheets.Item(6).Name = "Include ValueSets 5"
$wb.Worksheets.Item(7).Name = "Include ValueSets 6"
